$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.ClearFormats()
}

Set-TextValue $ws.Cells.Item(2, 4) "69.809.17"
Set-TextValue $ws.Cells.Item(2, 5) "  +5.17%  "

Set-TextValue $ws.Cells.Item(3, 4) "3.583.82"
Set-TextValue $ws.Cells.Item(3, 5) "  +4.95%  "

Set-TextValue $ws.Cells.Item(4, 5) "  +0.09%  "

Set-TextValue $ws.Cells.Item(5, 4) "586.13"
Set-TextValue $ws.Cells.Item(5, 5) "  +3.10%  "

Set-TextValue $ws.Cells.Item(6, 4) "189.88"
Set-TextValue $ws.Cells.Item(6, 5) "  +4.38%  "

Set-TextValue $ws.Cells.Item(7, 4) "0.644"
Set-TextValue $ws.Cells.Item(7, 5) "  +1.82%  "

Set-TextValue $ws.Cells.Item(8, 4) "3.578.79"
Set-TextValue $ws.Cells.Item(8, 5) "  +4.93%  "

Set-TextValue $ws.Cells.Item(9, 5) "  -0.06%  "

Set-TextValue $ws.Cells.Item(10, 5) "  -1.43%  "

Set-TextValue $ws.Cells.Item(11, 4) "0.657"
Set-TextValue $ws.Cells.Item(11, 5) "  +2.50%  "

Set-TextValue $ws.Cells.Item(12, 4) "57.62"
Set-TextValue $ws.Cells.Item(12, 5) "  +5.07%  "

Set-TextValue $ws.Cells.Item(13, 4) "0.0000290"
Set-TextValue $ws.Cells.Item(13, 5) "  +3.18%  "

Set-TextValue $ws.Cells.Item(14, 4) "9.71"
Set-TextValue $ws.Cells.Item(14, 5) "  +3.74%  "

Set-TextValue $ws.Cells.Item(15, 4) "4.156.17"
Set-TextValue $ws.Cells.Item(15, 5) "  +4.68%  "

Set-TextValue $ws.Cells.Item(18, 4) "69.823.87"
Set-TextValue $ws.Cells.Item(18, 5) "  +5.40%  "

Set-TextValue $ws.Cells.Item(19, 4) "12.45"
Set-TextValue $ws.Cells.Item(19, 5) "  +3.74%  "

Set-TextValue $ws.Cells.Item(20, 4) "0.120"
Set-TextValue $ws.Cells.Item(20, 5) "  +0.25%  "

Set-TextValue $ws.Cells.Item(21, 4) "1.04"
Set-TextValue $ws.Cells.Item(21, 5) "  +3.73%  "

Set-TextValue $ws.Cells.Item(22, 4) "493.94"
Set-TextValue $ws.Cells.Item(22, 5) "  +5.80%  "

Set-TextValue $ws.Cells.Item(23, 4) "17.58"
Set-TextValue $ws.Cells.Item(23, 5) "  +20.01%  "

Set-TextValue $ws.Cells.Item(24, 5) "  +7.30%  "

Set-TextValue $ws.Cells.Item(25, 4) "4.43"
Set-TextValue $ws.Cells.Item(25, 5) "  +6.85%  "

Set-TextValue $ws.Cells.Item(26, 4) "90.24"
Set-TextValue $ws.Cells.Item(26, 5) "  +0.55%  "

Set-TextValue $ws.Cells.Item(27, 5) "  +5.21%  "

Set-TextValue $ws.Cells.Item(28, 4) "11.06"
Set-TextValue $ws.Cells.Item(28, 5) "  +1.91%  "

Set-TextValue $ws.Cells.Item(29, 4) "9.36"
Set-TextValue $ws.Cells.Item(29, 5) "  +5.47%  "

Set-TextValue $ws.Cells.Item(30, 4) "32.10"
Set-TextValue $ws.Cells.Item(30, 5) "  +2.29%  "

Set-TextValue $ws.Cells.Item(31, 4) "7.52"
Set-TextValue $ws.Cells.Item(31, 5) "  +8.50%  "

Set-TextValue $ws.Cells.Item(34, 4) "0.116"
Set-TextValue $ws.Cells.Item(34, 5) "  +6.69%  "

Set-TextValue $ws.Cells.Item(35, 4) "65.16"
Set-TextValue $ws.Cells.Item(35, 5) "  +4.22%  "

Set-TextValue $ws.Cells.Item(36, 4) "0.0₃0811"
Set-TextValue $ws.Cells.Item(36, 5) "  +6.85%  "

Set-TextValue $ws.Cells.Item(37, 4) "0.404"
Set-TextValue $ws.Cells.Item(37, 5) "  +5.09%  "

Set-TextValue $ws.Cells.Item(38, 5) "  +0.09%  "

Set-TextValue $ws.Cells.Item(41, 4) "3.62"
Set-TextValue $ws.Cells.Item(41, 5) "  +1.07%  "

Set-TextValue $ws.Cells.Item(42, 4) "3.312.35"

Set-TextValue $ws.Cells.Item(43, 4) "3.06"
Set-TextValue $ws.Cells.Item(43, 5) "  +3.92%  "

Set-TextValue $ws.Cells.Item(44, 4) "0.0442"
Set-TextValue $ws.Cells.Item(44, 5) "  +3.78%  "

Set-TextValue $ws.Cells.Item(45, 4) "2.65"
Set-TextValue $ws.Cells.Item(45, 5) "  +4.68%  "

Set-TextValue $ws.Cells.Item(46, 5) "  +4.65%  "

Set-TextValue $ws.Cells.Item(47, 5) "  +1.49%  "

Set-TextValue $ws.Cells.Item(48, 4) "9.04"
Set-TextValue $ws.Cells.Item(48, 5) "  +5.17%  "

Set-TextValue $ws.Cells.Item(51, 4) "1.00"
Set-TextValue $ws.Cells.Item(51, 5) "  +0.28%  "

Set-TextValue $ws.Cells.Item(16, 2) "Chainlink"
Set-TextValue $ws.Cells.Item(16, 3) "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Cells.Item(16, 4) "19.28"
Set-TextValue $ws.Cells.Item(16, 5) "  +5.06%  "

Set-TextValue $ws.Cells.Item(17, 2) "WrappedEther"
Set-TextValue $ws.Cells.Item(17, 3) "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Cells.Item(17, 4) "3.587.18"
Set-TextValue $ws.Cells.Item(17, 5) "  +5.36%  "

Set-TextValue $ws.Cells.Item(32, 2) "Cosmos"
Set-TextValue $ws.Cells.Item(32, 3) "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Cells.Item(32, 4) "12.17"
Set-TextValue $ws.Cells.Item(32, 5) "  +5.11%  "

Set-TextValue $ws.Cells.Item(33, 2) "Bittensor"
Set-TextValue $ws.Cells.Item(33, 3) "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Cells.Item(33, 4) "614.56"
Set-TextValue $ws.Cells.Item(33, 5) "  +5.01%  "

Set-TextValue $ws.Cells.Item(39, 2) "Kaspa"
Set-TextValue $ws.Cells.Item(39, 3) "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Cells.Item(39, 4) "0.146"
Set-TextValue $ws.Cells.Item(39, 5) "  +0.70%  "

Set-TextValue $ws.Cells.Item(40, 2) "InjectiveProtocol"
Set-TextValue $ws.Cells.Item(40, 3) "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Cells.Item(40, 4) "37.81"
Set-TextValue $ws.Cells.Item(40, 5) "  +3.67%  "

Set-TextValue $ws.Cells.Item(49, 2) "LidoDAOToken"
Set-TextValue $ws.Cells.Item(49, 3) "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Cells.Item(49, 4) "3.32"
Set-TextValue $ws.Cells.Item(49, 5) "  +5.65%  "

Set-TextValue $ws.Cells.Item(50, 2) "dogwifhat"
Set-TextValue $ws.Cells.Item(50, 3) "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Cells.Item(50, 4) "2.68"
Set-TextValue $ws.Cells.Item(50, 5) "  -4.24%  "

Write-Host "Applied cryptos update"